$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 102.1928919321432
$ws.Range("B2").Value = 0.1302852986657663
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 125.1557687919098
$ws.Range("B3").Value = -0.06514073549644994
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 96.93063574909577
$ws.Range("B4").Value = 0.1750698377668048
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 108.2889579438499
$ws.Range("B5").Value = 0.1036522059779796
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = 106.6518651426058
$ws.Range("B6").Value = 0.1172030291539252
$ws.Range("C6").Value = 5

$wb.Save()
